# This script applies the "added harvard case classification" update.
# Adding the new Harvard doctor-case classification shifts the existing
# "average_doctor" figures into a new "average_doctor_old" column (BQ) and
# populates the "average_doctor" column (BP) with freshly recomputed values,
# while several of the underlying per-app stat cells in rows 4-13 change too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Swap header labels for columns BP/BQ.
$ws.Range("BP1").Value2 = "average_doctor_old"
$ws.Range("BQ1").Value2 = "average_doctor"

# Updated stats for row 4
$ws.Range("E4").Value2 = 0.481
$ws.Range("F4").Value2 = 0.053
$ws.Range("G4").Value2 = 0.23
$ws.Range("N4").Value2 = 0.485
$ws.Range("O4").Value2 = 0.057
$ws.Range("P4").Value2 = 0.24
$ws.Range("Q4").Value2 = 0.052
$ws.Range("R4").Value2 = 0.035
$ws.Range("S4").Value2 = 0.188
$ws.Range("W4").Value2 = 0.367
$ws.Range("AI4").Value2 = 0.403
$ws.Range("AJ4").Value2 = 0.096
$ws.Range("AK4").Value2 = 0.309
$ws.Range("AU4").Value2 = 0.24
$ws.Range("AW4").Value2 = 0.163
$ws.Range("BA4").Value2 = 2.041
$ws.Range("BB4").Value2 = 0.142
$ws.Range("BC4").Value2 = 0.377
$ws.Range("BG4").Value2 = 0.722
$ws.Range("BH4").Value2 = 0.142
$ws.Range("BI4").Value2 = 0.377
$ws.Range("BM4").Value2 = 0.75
$ws.Range("BN4").Value2 = 0.064
$ws.Range("BO4").Value2 = 0.253
$ws.Range("BP4").Value2 = 0.68
$ws.Range("BQ4").Value2 = 0.761

# Updated stats for row 5
$ws.Range("E5").Value2 = 0.606
$ws.Range("F5").Value2 = 0.058
$ws.Range("G5").Value2 = 0.241
$ws.Range("N5").Value2 = 0.733
$ws.Range("O5").Value2 = 0.065
$ws.Range("P5").Value2 = 0.255
$ws.Range("Q5").Value2 = 0.035
$ws.Range("R5").Value2 = 0.014
$ws.Range("S5").Value2 = 0.12
$ws.Range("W5").Value2 = 0.335
$ws.Range("X5").Value2 = 0.1
$ws.Range("Y5").Value2 = 0.316
$ws.Range("AI5").Value2 = 0.404
$ws.Range("AJ5").Value2 = 0.092
$ws.Range("AK5").Value2 = 0.303
$ws.Range("AU5").Value2 = 0.448
$ws.Range("AV5").Value2 = 0.079
$ws.Range("AW5").Value2 = 0.28
$ws.Range("BA5").Value2 = 1.306
$ws.Range("BB5").Value2 = 0.074
$ws.Range("BC5").Value2 = 0.271
$ws.Range("BG5").Value2 = 0.383
$ws.Range("BH5").Value2 = 0.051
$ws.Range("BI5").Value2 = 0.225
$ws.Range("BM5").Value2 = 0.525
$ws.Range("BN5").Value2 = 0.047
$ws.Range("BO5").Value2 = 0.217
$ws.Range("BP5").Value2 = 0.435
$ws.Range("BQ5").Value2 = 0.458

# Updated stats for row 6
$ws.Range("E6").Value2 = 0.536
$ws.Range("N6").Value2 = 0.584
$ws.Range("Q6").Value2 = 0.042
$ws.Range("W6").Value2 = 0.35
$ws.Range("AI6").Value2 = 0.403
$ws.Range("AU6").Value2 = 0.313
$ws.Range("BA6").Value2 = 1.586
$ws.Range("BG6").Value2 = 0.5
$ws.Range("BM6").Value2 = 0.618
$ws.Range("BP6").Value2 = 0.529
$ws.Range("BQ6").Value2 = 0.569

# Updated stats for row 7
$ws.Range("E7").Value2 = 0.576
$ws.Range("N7").Value2 = 0.665
$ws.Range("Q7").Value2 = 0.037
$ws.Range("W7").Value2 = 0.341
$ws.Range("AI7").Value2 = 0.404
$ws.Range("AU7").Value2 = 0.382
$ws.Range("BA7").Value2 = 1.405
$ws.Range("BG7").Value2 = 0.423
$ws.Range("BM7").Value2 = 0.5590000000000001
$ws.Range("BP7").Value2 = 0.468
$ws.Range("BQ7").Value2 = 0.497

# Updated stats for row 8
$ws.Range("E8").Value2 = 0.705
$ws.Range("F8").Value2 = 0.07199999999999999
$ws.Range("G8").Value2 = 0.268
$ws.Range("N8").Value2 = 0.821
$ws.Range("O8").Value2 = 0.044
$ws.Range("P8").Value2 = 0.21
$ws.Range("Q8").Value2 = 0.038
$ws.Range("W8").Value2 = 0.407
$ws.Range("X8").Value2 = 0.122
$ws.Range("Y8").Value2 = 0.349
$ws.Range("AI8").Value2 = 0.472
$ws.Range("AJ8").Value2 = 0.14
$ws.Range("AK8").Value2 = 0.374
$ws.Range("AU8").Value2 = 0.392
$ws.Range("AV8").Value2 = 0.08500000000000001
$ws.Range("AW8").Value2 = 0.292
$ws.Range("BA8").Value2 = 1.77
$ws.Range("BB8").Value2 = 0.109
$ws.Range("BC8").Value2 = 0.33
$ws.Range("BG8").Value2 = 0.57
$ws.Range("BH8").Value2 = 0.11
$ws.Range("BI8").Value2 = 0.331
$ws.Range("BM8").Value2 = 0.673
$ws.Range("BN8").Value2 = 0.062
$ws.Range("BO8").Value2 = 0.249
$ws.Range("BP8").Value2 = 0.59
$ws.Range("BQ8").Value2 = 0.626

# Updated stats for row 9
$ws.Range("E9").Value2 = 0.667
$ws.Range("F9").Value2 = 0.222
$ws.Range("G9").Value2 = 0.471
$ws.Range("N9").Value2 = 0.762
$ws.Range("O9").Value2 = 0.181
$ws.Range("P9").Value2 = 0.426
$ws.Range("W9").Value2 = 0.31
$ws.Range("X9").Value2 = 0.214
$ws.Range("Y9").Value2 = 0.462
$ws.Range("AI9").Value2 = 0.429
$ws.Range("AJ9").Value2 = 0.245
$ws.Range("AK9").Value2 = 0.495
$ws.Range("BA9").Value2 = 1.738
$ws.Range("BB9").Value2 = 0.249
$ws.Range("BC9").Value2 = 0.499
$ws.Range("BG9").Value2 = 0.619
$ws.Range("BH9").Value2 = 0.236
$ws.Range("BI9").Value2 = 0.486
$ws.Range("BM9").Value2 = 0.643
$ws.Range("BN9").Value2 = 0.23
$ws.Range("BO9").Value2 = 0.479
$ws.Range("BP9").Value2 = 0.579
$ws.Range("BQ9").Value2 = 0.618

# Updated stats for row 10
$ws.Range("E10").Value2 = 0.8100000000000001
$ws.Range("F10").Value2 = 0.154
$ws.Range("G10").Value2 = 0.393
$ws.Range("N10").Value2 = 0.952
$ws.Range("O10").Value2 = 0.045
$ws.Range("P10").Value2 = 0.213
$ws.Range("W10").Value2 = 0.524
$ws.Range("X10").Value2 = 0.249
$ws.Range("Y10").Value2 = 0.499
$ws.Range("AI10").Value2 = 0.5
$ws.Range("AJ10").Value2 = 0.25
$ws.Range("AK10").Value2 = 0.5
$ws.Range("AU10").Value2 = 0.381
$ws.Range("AV10").Value2 = 0.236
$ws.Range("AW10").Value2 = 0.486
$ws.Range("BA10").Value2 = 2.19
$ws.Range("BB10").Value2 = 0.214
$ws.Range("BC10").Value2 = 0.462
$ws.Range("BG10").Value2 = 0.6899999999999999
$ws.Range("BH10").Value2 = 0.214
$ws.Range("BI10").Value2 = 0.462
$ws.Range("BM10").Value2 = 0.8100000000000001
$ws.Range("BN10").Value2 = 0.154
$ws.Range("BO10").Value2 = 0.393
$ws.Range("BP10").Value2 = 0.73
$ws.Range("BQ10").Value2 = 0.764

# Updated stats for row 11
$ws.Range("E11").Value2 = 0.857
$ws.Range("F11").Value2 = 0.122
$ws.Range("G11").Value2 = 0.35
$ws.Range("N11").Value2 = 0.952
$ws.Range("O11").Value2 = 0.045
$ws.Range("P11").Value2 = 0.213
$ws.Range("W11").Value2 = 0.524
$ws.Range("X11").Value2 = 0.249
$ws.Range("Y11").Value2 = 0.499
$ws.Range("AI11").Value2 = 0.571
$ws.Range("AJ11").Value2 = 0.245
$ws.Range("AK11").Value2 = 0.495
$ws.Range("AU11").Value2 = 0.548
$ws.Range("AV11").Value2 = 0.248
$ws.Range("AW11").Value2 = 0.498
$ws.Range("BA11").Value2 = 2.19
$ws.Range("BB11").Value2 = 0.214
$ws.Range("BC11").Value2 = 0.462
$ws.Range("BG11").Value2 = 0.6899999999999999
$ws.Range("BH11").Value2 = 0.214
$ws.Range("BI11").Value2 = 0.462
$ws.Range("BM11").Value2 = 0.8100000000000001
$ws.Range("BN11").Value2 = 0.154
$ws.Range("BO11").Value2 = 0.393
$ws.Range("BP11").Value2 = 0.73
$ws.Range("BQ11").Value2 = 0.77

# Updated stats for row 12
$ws.Range("E12").Value2 = 1.417
$ws.Range("F12").Value2 = 0.854
$ws.Range("G12").Value2 = 0.924
$ws.Range("N12").Value2 = 1.25
$ws.Range("O12").Value2 = 0.287
$ws.Range("P12").Value2 = 0.536
$ws.Range("W12").Value2 = 1.5
$ws.Range("X12").Value2 = 0.432
$ws.Range("Y12").Value2 = 0.657
$ws.Range("AI12").Value2 = 1.583
$ws.Range("AJ12").Value2 = 1.493
$ws.Range("AK12").Value2 = 1.222
$ws.Range("AU12").Value2 = 2.88
$ws.Range("AV12").Value2 = 3.466
$ws.Range("AW12").Value2 = 1.862
$ws.Range("BA12").Value2 = 3.787
$ws.Range("BB12").Value2 = 0.449
$ws.Range("BC12").Value2 = 0.67
$ws.Range("BG12").Value2 = 1.138
$ws.Range("BH12").Value2 = 0.188
$ws.Range("BI12").Value2 = 0.433
$ws.Range("BM12").Value2 = 1.235
$ws.Range("BN12").Value2 = 0.239
$ws.Range("BO12").Value2 = 0.489
$ws.Range("BP12").Value2 = 1.262
$ws.Range("BQ12").Value2 = 1.248

# Updated stats for row 13
$ws.Range("E13").Value2 = 1.415
$ws.Range("F13").Value2 = 0.295
$ws.Range("G13").Value2 = 0.543
$ws.Range("N13").Value2 = 1.737
$ws.Range("O13").Value2 = 0.466
$ws.Range("P13").Value2 = 0.6830000000000001
$ws.Range("W13").Value2 = 0.985
$ws.Range("X13").Value2 = 0.199
$ws.Range("Y13").Value2 = 0.446
$ws.Range("AI13").Value2 = 1.159
$ws.Range("AJ13").Value2 = 0.312
$ws.Range("AK13").Value2 = 0.5580000000000001
$ws.Range("AU13").Value2 = 2.048
$ws.Range("AV13").Value2 = 0.344
$ws.Range("AW13").Value2 = 0.587
$ws.Range("BA13").Value2 = 2.187
$ws.Range("BB13").Value2 = 0.278
$ws.Range("BC13").Value2 = 0.527
$ws.Range("BG13").Value2 = 0.547
$ws.Range("BH13").Value2 = 0.05
$ws.Range("BI13").Value2 = 0.224
$ws.Range("BM13").Value2 = 0.787
$ws.Range("BN13").Value2 = 0.163
$ws.Range("BO13").Value2 = 0.403
$ws.Range("BP13").Value2 = 0.729
$ws.Range("BQ13").Value2 = 0.667
